# edit.ps1
# Applies the "updated documentation for project" changes to the
# Progress Journal document.

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Title line: merge the " " run and "Progress Journal" run into
#    a single run reading " Progress Journal".
#    (Find/Replace over the combined text collapses the two runs
#    that get matched into one run bearing the replacement text.)
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" Progress Journal", $false, $false, $false, $false, $false, `
                   $true, 1, $false, " Progress Journal", 2) | Out-Null

# ---------------------------------------------------------------
# 2) "Task Progress" heading: merge the "T" run and "ask Progress"
#    run into a single run reading "Task Progress".
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Task Progress", $false, $false, $false, $false, $false, `
                   $true, 1, $false, "Task Progress", 2) | Out-Null

# ---------------------------------------------------------------
# 3) Week 4 entry: append the week's write-up after "Week 4 - ".
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Week 4 " + [char]0x2013 + " ") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("Didn" + [char]0x2019 + "t achieve lots this week as Henry was sick with covid, I had to work late Monday night and missed class. I tried to do some stuff Wednesday setting up collectible resources for the game but ran into version control issues again. The project didn" + [char]0x2019 + "t appear to upload correctly to git and when I tried to work on a version of it. I encountered errors that I couldn" + [char]0x2019 + "t fix, and the project was a bit broken. I continued to work on some of the documentation instead.") | Out-Null

# ---------------------------------------------------------------
# 4) Week 5 entry: append the week's write-up after "Week 5 - ".
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Week 5 " + [char]0x2013 + " ") | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("In Henry" + [char]0x2019 + "s time off with covid he started a new version of the game and brought that on Monday for us to continue with. The game was the same concept, and Henry had also implemented resources and the base for the progression system. We spent a lot of time trying to upload the new project to git which was again being difficult due to Git 100mb file restriction and GitLFS being a pain to set up. We changed to bit bucket and finally got the project uploaded successfully so we can work on it collaboratively. Came up with some more ideas to make the game better regarding progression, enemies and a minimap.") | Out-Null

# ---------------------------------------------------------------
# 5) Issues Found list: append three new bullet items after the
#    "... a bit more basic." item, re-using the same ListParagraph /
#    numbered-list formatting as the paragraph before them.
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("if it proves to be too difficult then we will just get something a bit more basic.") | Out-Null
$rng.Collapse(0)

$rng.InsertParagraphAfter() | Out-Null
$rng.Collapse(0)
$rng.MoveStart(1, 1) | Out-Null
$rng.InsertAfter("GitHub 100mb file limit proved too difficult to upload the projects without breaking changed to Bit Bucket.") | Out-Null

$rng.Collapse(0)
$rng.InsertParagraphAfter() | Out-Null
$rng.Collapse(0)
$rng.MoveStart(1, 1) | Out-Null
$rng.InsertAfter("Bit Bucket account was painful to set up and clone a repository due to some weird access protocol they implemented. Nick helped us figure out how to set it up correctly and get the version control working.") | Out-Null

$rng.Collapse(0)
$rng.InsertParagraphAfter() | Out-Null
$rng.Collapse(0)
$rng.MoveStart(1, 1) | Out-Null
$rng.InsertAfter("The new version of the project that Henry created was in a newer version of unity that wasn" + [char]0x2019 + "t installed on the AIE computers which slowed our development for the night. ") | Out-Null
